# Updated_Format_1 (Text Wrapper+Coloana B mai mica)
# Applies:
#   1. A3 label gets two leading spaces ("  Numele și Prenumele:")
#   2. Light-gray "counter" font (fontId 1) becomes much lighter (E3E3E3) and
#      is centered (horizontal+vertical) for every cell that uses it
#   3. Column B is narrowed (was ~250 chars wide, now ~200 chars wide)
#   4. Every "content" fill style gets WrapText turned on (and the orange
#      header-label style additionally becomes horizontally centered)
#   5. Each sentence block gets a word/character count dropped into column A
#      of its first filler row, styled like the other "counter" cells

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# ---------------------------------------------------------------------
# 1) Update the label text in A3 to have a leading two-space indent
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "  Numele și Prenumele:"

# ---------------------------------------------------------------------
# 2) Re-color + center the "counter" style (fontId 1 / s=2) cells.
#    Apply to A1 first (creates the new style), then copy that
#    formatting onto the rest so they all reuse the same style record.
# ---------------------------------------------------------------------
$ws.Range("A1").Font.Color = 14935011   # 0xE3E3E3
$ws.Range("A1").HorizontalAlignment = $xlCenter
$ws.Range("A1").VerticalAlignment = $xlCenter

$ws.Range("A1").Copy() | Out-Null
$counterCells = @("A6","A18","A30","A42","A54","A66","A78","A90","A102","A114","A126")
foreach ($addr in $counterCells) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Narrow column B
# ---------------------------------------------------------------------
$ws.Range("B1").ColumnWidth = 199.833333

# ---------------------------------------------------------------------
# 4) Turn on WrapText for every content-fill style; the orange
#    "header label" style (fillId 5 / s=6) additionally gets centered.
# ---------------------------------------------------------------------

# s=6 : A5, A17, A29, A41, A53, A65, A77, A89, A101, A113, A125 (orange header label, col A)
$ws.Range("A5").HorizontalAlignment = $xlCenter
$ws.Range("A5").VerticalAlignment = $xlCenter
$ws.Range("A5").WrapText = $true
$ws.Range("A5").Copy() | Out-Null
$s6cells = @("A17","A29","A41","A53","A65","A77","A89","A101","A113","A125")
foreach ($addr in $s6cells) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# s=7 : B5, B17, B29, B41, B53, B65, B77, B89, B101, B113, B125 (blue header text, col B)
$ws.Range("B5").WrapText = $true
$ws.Range("B5").Copy() | Out-Null
$s7cells = @("B17","B29","B41","B53","B65","B77","B89","B101","B113","B125")
foreach ($addr in $s7cells) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# s=8 : every other blank filler cell (red tint)
$ws.Range("B6").WrapText = $true
$ws.Range("B6").Copy() | Out-Null
$s8cells = @("B8","B10","B12","B14","B20","B22","B24","B26","B32","B34","B36","B38","B42","B44","B46","B48","B50","B56","B58","B60","B62","B66","B68","B70","B72","B74","B78","B80","B82","B84","B86","B90","B92","B94","B96","B98","B102","B104","B106","B108","B110","B114","B116","B118","B120","B122","B126","B128","B130","B132","B134")
foreach ($addr in $s8cells) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# s=9 : alternating blank filler cell (red, slightly different tint)
$ws.Range("B7").WrapText = $true
$ws.Range("B7").Copy() | Out-Null
$s9cells = @("B9","B11","B13","B15","B19","B21","B23","B25","B27","B31","B33","B35","B37","B39","B43","B45","B47","B49","B51","B57","B59","B61","B63","B67","B69","B71","B73","B75","B79","B81","B83","B85","B87","B91","B93","B95","B97","B99","B103","B105","B107","B109","B111","B115","B117","B119","B121","B123","B127","B129","B131","B133","B135")
foreach ($addr in $s9cells) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# s=10 : example-sentence cells (light blue)
$ws.Range("B18").WrapText = $true
$ws.Range("B18").Copy() | Out-Null
$s10cells = @("B30","B54")
foreach ($addr in $s10cells) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# s=11 : the single extra example-sentence cell (green)
$ws.Range("B55").WrapText = $true

# ---------------------------------------------------------------------
# 5) Drop the sentence-count numbers into column A of the first filler
#    row of each block, matching the "counter" style used in A1/A6/...
# ---------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null

$ws.Range("A7").Value = 10
$ws.Range("A19").Value = 9
$ws.Range("A31").Value = 9
$ws.Range("A43").Value = 10
$ws.Range("A55").Value = 8
$ws.Range("A67").Value = 10
$ws.Range("A79").Value = 10
$ws.Range("A91").Value = 10
$ws.Range("A103").Value = 10
$ws.Range("A115").Value = 10
$ws.Range("A127").Value = 10

$countCells = @("A7","A19","A31","A43","A55","A67","A79","A91","A103","A115","A127")
foreach ($addr in $countCells) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0
